# audiences-process.pptx edit:
#  1. Bump the "datetimeFigureOut" date placeholder text (01/02/2023 -> 02/08/2023)
#     on the slide master and on every slide layout.
#  2. Widen the "TextBox 11" shape on slide 1 and extend its label text to
#     mention "journeys" as well as "campaigns".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder text, master + all layouts
# ---------------------------------------------------------------------------
$newDate = "02/08/2023"

function Update-DatePlaceholder($shapes, $newText) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shape = $shapes.Item($j)
        if ($shape.Name -like "Date Placeholder*") {
            $shape.TextFrame.TextRange.Text = $newText
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes $newDate

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Update-DatePlaceholder $layout.Shapes $newDate
}

# ---------------------------------------------------------------------------
# 2) "TextBox 11" on slide 1: wider box + updated copy
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes
for ($k = 1; $k -le $shapes.Count; $k++) {
    $shp = $shapes.Item($k)
    if ($shp.Name -eq "TextBox 11") {
        $shp.TextFrame.TextRange.Text = "Target audience(s) in campaigns and journeys"
        [double]$targetWidthPt = 2238215 / 12700.0
        $shp.Width = $targetWidthPt
    }
}
